$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row renames ---
$ws.Range("G1").Value = "nextLev"
$ws.Range("J1").Value = "atkType"

# --- Row 2: 加农炮1 ---
$ws.Range("H2").Value = "TowerImg/1_1"
$ws.Range("K2").Value = "eff/Fire1"

# --- Row 3: 加农炮2 ---
$ws.Range("H3").Value = "TowerImg/1_2"
$ws.Range("K3").Value = "eff/Fire1"

# --- Row 4: 加农炮3 ---
$ws.Range("H4").Value = "TowerImg/1_3"
$ws.Range("K4").Value = "eff/Fire1"

# --- Row 5: 机枪炮1 ---
$ws.Range("H5").Value = "TowerImg/2_1"
$ws.Range("K5").Value = "eff/Fire1"

# --- Row 6: 机枪炮2 ---
$ws.Range("H6").Value = "TowerImg/2_2"
$ws.Range("K6").Value = "eff/Fire1"

# --- Row 7: 机枪炮3 ---
$ws.Range("H7").Value = "TowerImg/2_3"
$ws.Range("K7").Value = "eff/Fire1"

# --- Row 8: 闪电炮1 -> 魔法炮1 ---
$ws.Range("B8").Value = "魔法炮1"
$ws.Range("H8").Value = "TowerImg/3_1"
$ws.Range("K8").Value = "eff/Fire2"

# --- Row 9: 闪电炮2 -> 魔法炮2 ---
$ws.Range("B9").Value = "魔法炮2"
$ws.Range("H9").Value = "TowerImg/3_2"
$ws.Range("K9").Value = "eff/Fire2"

# --- Row 10: 闪电炮3 -> 魔法炮3 ---
$ws.Range("B10").Value = "魔法炮3"
$ws.Range("H10").Value = "TowerImg/3_3"
$ws.Range("K10").Value = "eff/Fire2"

# --- Column widths ---
# (column B's custom width is removed in the target - leave it untouched/default)
$ws.Columns.Item(5).ColumnWidth = 10
$ws.Columns.Item(6).ColumnWidth = 12.142857142857142
$ws.Columns.Item(8).ColumnWidth = 14.142857142857142
$ws.Columns.Item(9).ColumnWidth = 14.285714285714286
$ws.Columns.Item(11).ColumnWidth = 16.857142857142858

# --- Selection ---
$ws.Range("A8").Select()
